# DSA Binary Search Answers Practice
# Applies:
#  1. Fix mislabeled Difficulty on rows 39/40 ("E" -> "Easy")
#  2. Add a bold, centered grand-total formula in H1 (=SUM(H3:H5))
#  3. Append two new LeetCode "Medium" binary-search questions (rows 55/56)
#     with their hyperlinks
#  4. Update the view's zoom + selection to reflect the new data range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the two rows that used the stray single-letter "E" code instead
#        of the full "Easy" difficulty word already used everywhere else.
$ws.Range("D39").Value2 = "Easy"
$ws.Range("D40").Value2 = "Easy"

# --- 2. New bold/centered grand total cell above the Easy/Medium/Hard counts
$ws.Range("H1").Formula = "=SUM(H3:H5)"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4108

# --- 3. Append the two new rows, copying the formatting used by the last
#        existing data row (54) so borders/fonts match the rest of the table.
$ws.Range("A54:F54").Copy()
$ws.Range("A55:F55").PasteSpecial(-4122)
$ws.Range("A56:F56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill row 56 first so its brand-new strings are interned before row 55's,
# matching the original authoring order of the shared-string table.
$ws.Range("A56").Value2 = 54
$ws.Range("B56").Value2 = "1011. Capacity To Ship Packages Within D Days"
$ws.Range("C56").Value2 = "LeetCode"
$ws.Range("D56").Value2 = "Medium"
$ws.Range("E56").Value2 = "https://leetcode.com/problems/capacity-to-ship-packages-within-d-days/description/"
$ws.Range("F56").Value2 = "|"
$ws.Hyperlinks.Add($ws.Range("E56"), $ws.Range("E56").Value2) | Out-Null
# Hyperlinks.Add() stamps an extra (harmless) applyFont flag onto the cell's
# style, which would otherwise mint a redundant style entry; re-paste the
# original hyperlink-cell formatting on top so E56 keeps reusing style 2.
$ws.Range("E54").Copy()
$ws.Range("E56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A55").Value2 = 53
$ws.Range("B55").Value2 = "1283. Find the Smallest Divisor Given a Threshold"
$ws.Range("C55").Value2 = "LeetCode"
$ws.Range("D55").Value2 = "Medium"
$ws.Range("E55").Value2 = "https://leetcode.com/problems/find-the-smallest-divisor-given-a-threshold/description/"
$ws.Range("F55").Value2 = "|"
$ws.Hyperlinks.Add($ws.Range("E55"), $ws.Range("E55").Value2) | Out-Null
$ws.Range("E54").Copy()
$ws.Range("E55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Refresh the view: zoom to 90% and leave the selection on the new
#        last row, same as the author's session.
$excel.ActiveWindow.Zoom = 90
$ws.Range("B56").Select()
